# Include process sets splitting processes by sector
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sets-Proc")

# New rows 14-20: process sets split out by sector (AGR, SRV, IND, PWR, RSD, SUP, TRA)
$rows = @(
    @{ Row=14; A=$null;  B="A*,FT-AGR*";        F="PRC_AGR"; G="All AGR processes" },
    @{ Row=15; A=$null;  B="S-*,FT-SRV*";       F="PRC_SRV"; G="All SRV processes" },
    @{ Row=16; A="-IRE"; B="I*,FT-IND*";        F="PRC_IND"; G="All IND processes" },
    @{ Row=17; A=$null;  B="P*,FT-PWR*,*GRID*"; F="PRC_PWR"; G="All PWR processes" },
    @{ Row=18; A=$null;  B="R*,FT-RSD*";        F="PRC_RSD"; G="All RSD processes" },
    @{ Row=19; A=$null;  B="S*,FT-SUP*,-S-*";   F="PRC_SUP"; G="All SUP processes" },
    @{ Row=20; A="-IRE"; B="T*,FT-TRA*";        F="PRC_TRA"; G="All TRA processes" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    if ($r.A -ne $null) {
        $ws.Cells.Item($rowNum, 1).Value = $r.A
    }
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = "AND"
    $ws.Cells.Item($rowNum, 9).Value = "OR"
    $ws.Cells.Item($rowNum, 10).Value = "AND"
    $ws.Cells.Item($rowNum, 11).Value = "OR"
}

# SRV_Sets-Proc previously had focus (tabSelected + selection D18); move its
# selection but leave it no longer the active tab.
$srv = $wb.Worksheets.Item("SRV_Sets-Proc")
[void]$srv.Activate()
[void]$srv.Range("H31").Select()

# Sets-Proc becomes the active/selected tab (activeTab index 1), with
# selection moved to B27.
[void]$ws.Activate()
[void]$ws.Range("B27").Select()
